# Weekly update: insert two new price rows (Kiwi, Macroferia Regional de
# Talca) at the top of the data block that starts at row 338, pushing the
# existing rows (338:362) down by two positions (-> 340:364).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("338:339").Insert()

# New row 338: Kiwi, Hayward, calidad "Primera"
$ws.Cells.Item(338, 1).Value  = 5
$ws.Cells.Item(338, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(338, 3).Value  = "Maule"
$ws.Cells.Item(338, 4).Value  = 44826
$ws.Cells.Item(338, 5).Value  = 7
$ws.Cells.Item(338, 6).Value  = "Fruta"
$ws.Cells.Item(338, 7).Value  = 100101
$ws.Cells.Item(338, 8).Value  = "Berries"
$ws.Cells.Item(338, 9).Value  = 100101007
$ws.Cells.Item(338, 10).Value = "Kiwi"
$ws.Cells.Item(338, 11).Value = "Hayward"
$ws.Cells.Item(338, 12).Value = "Primera"
$ws.Cells.Item(338, 13).Value = 250
$ws.Cells.Item(338, 14).Value = 8000
$ws.Cells.Item(338, 15).Value = 8000
$ws.Cells.Item(338, 16).Value = 8000
$ws.Cells.Item(338, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(338, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(338, 19).Value = 444
$ws.Cells.Item(338, 20).Value = 18

# New row 339: Kiwi, Hayward, calidad "Segunda"
$ws.Cells.Item(339, 1).Value  = 5
$ws.Cells.Item(339, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(339, 3).Value  = "Maule"
$ws.Cells.Item(339, 4).Value  = 44826
$ws.Cells.Item(339, 5).Value  = 7
$ws.Cells.Item(339, 6).Value  = "Fruta"
$ws.Cells.Item(339, 7).Value  = 100101
$ws.Cells.Item(339, 8).Value  = "Berries"
$ws.Cells.Item(339, 9).Value  = 100101007
$ws.Cells.Item(339, 10).Value = "Kiwi"
$ws.Cells.Item(339, 11).Value = "Hayward"
$ws.Cells.Item(339, 12).Value = "Segunda"
$ws.Cells.Item(339, 13).Value = 200
$ws.Cells.Item(339, 14).Value = 6000
$ws.Cells.Item(339, 15).Value = 6000
$ws.Cells.Item(339, 16).Value = 6000
$ws.Cells.Item(339, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(339, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(339, 19).Value = 333
$ws.Cells.Item(339, 20).Value = 18
